# Weekly update: a new week's worth of "Cebolla" (onion) price data for
# Femacal de La Calera (Coquimbo) is inserted at the top of the existing
# date-ordered data block (rows 621-636), pushing the prior rows down by
# five rows (621->626 ... 636->641) and growing the sheet from R636 to R641.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows at 621, shifting existing rows 621:636 down to 626:641.
$ws.Rows("621:625").Insert()

# New week's data (Fecha serial 44448) for the five "Calidad" buckets that
# appear for every weekly block in this sheet.
$newRows = @(
    @{ Row=621; H="Morada(o)";        I="Primera";     J=130;   K=8500; L=9000; M=8769; N="$/malla 18 kilos";                             O="Perú";                  P=487; Q=18 },
    @{ Row=622; H="Sin especificar";  I="1a (guarda)"; J=185;   K=4500; L=5000; M=4757; N="$/malla 18 kilos";                             O="Provincia de Quillota"; P=264; Q=18 },
    @{ Row=623; H="Sin especificar";  I="1a nueva(o)"; J=11300; K=4300; L=4500; M=4403; N="$/paquete 20 unidades (volumen en unidades)";  O="Provincia de Quillota"; P=220; Q=20 },
    @{ Row=624; H="Sin especificar";  I="2a (guarda)"; J=80;    K=4000; L=4000; M=4000; N="$/malla 18 kilos";                             O="Provincia de Quillota"; P=222; Q=18 },
    @{ Row=625; H="Sin especificar";  I="2a nueva(o)"; J=5800;  K=3500; L=3500; M=3500; N="$/paquete 20 unidades (volumen en unidades)";  O="Provincia de Quillota"; P=175; Q=20 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 3
    $ws.Cells.Item($row, 2).Value = "Femacal de La Calera"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = 44448
    $ws.Cells.Item($row, 5).Value = 5
    $ws.Cells.Item($row, 6).Value = 100112004
    $ws.Cells.Item($row, 7).Value = "Cebolla"
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
